$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

$ws.Range("B5").Copy()
$ws.Range("B6:B7").PasteSpecial(-4122)

$ws.Range("B6").Value = 43332
$ws.Range("C6").Value = 124

$ws.Range("B7").Value = 43333
$ws.Range("C7").Value = 126

$ws.Range("C8").Select()
